$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Database" stage tasks (rows 4-9, column F = Status) as done (TRUE)
$ws.Range("F4:F9").Value = $true

# Update the active cell selection to match the latest edit position
$ws.Activate()
$ws.Range("J7").Select()
